$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8771
$ws1.Range("F3").Value = 94
$ws1.Range("F4").Value = 236
$ws1.Range("F5").Value = 99
$ws1.Range("F6").Value = 1421
$ws1.Range("F7").Value = 1386
$ws1.Range("F8").Value = 240
$ws1.Range("F9").Value = 39
$ws1.Range("F10").Value = 284
$ws1.Range("F11").Value = 82

# Sheet "全部类型" (All Types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8771
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 236
$ws4.Range("F5").Value = 99
$ws4.Range("F6").Value = 1421
$ws4.Range("F7").Value = 1386
$ws4.Range("F8").Value = 240
$ws4.Range("F10").Value = 39
$ws4.Range("F11").Value = 284
$ws4.Range("F12").Value = 82
